$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: replace the "EEUFC1V331 (Cout)" part with the new "ECA-1VM331" component,
# update its unit price, and point the product link to the new datasheet.
$ws.Range("A11").Value = "ECA-1VM331"
$ws.Range("B11").Value = 0.41
$ws.Range("F11").Value = "http://goo.gl/3D2k0X"

# Update the sheet view: drop the frozen top-left cell, zoom out to 70%,
# and move the active selection to C12.
$ws.Application.ActiveWindow.Zoom = 70
$ws.Range("C12").Select()
